# Utilities.xlsx edit: "Set up @timing framework to track each step."
#
# Inserts a new worksheet "@timing framework" as the first sheet in the
# workbook, containing timing-delta data for @timing[pt_econfig] vs
# @time_in_e_config, and nudges a couple of leftover view-state bits
# (active selections) on two of the existing sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new sheet at the very front of the workbook.
#    Worksheets.Add() with no args inserts before the currently active
#    sheet, which (on a freshly-opened workbook) is sheet 1.
# ---------------------------------------------------------------------
$timing = $wb.Worksheets.Add([System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value)
$timing.Name = "@timing framework"

# Header row
$timing.Cells.Item(1,1).Value = "Run"
$timing.Cells.Item(1,2).Formula = "'@timing[pt_econfig]"
$timing.Cells.Item(1,3).Formula = "'@time_in_e_config"
$timing.Cells.Item(1,4).Value = "Delta"

# Copy the quote-prefixed text formatting from the B1/C1 headers down
# into the (otherwise blank) summary-label cells below them, matching
# the look of the other timing-comparison sheets in this workbook.
$timing.Range("B1:C1").Copy() | Out-Null
$timing.Range("B2:C4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Summary formulas (rows 2-4): average / stdev / relative stdev of the
# per-run deltas in column D.
$timing.Cells.Item(2,4).Formula = "=AVERAGE(D5:D31)"
$timing.Cells.Item(3,4).Formula = "=STDEVA(D5:D31)"
$timing.Cells.Item(4,4).Formula = "=D3/D2"
$timing.Cells.Item(4,4).NumberFormat = "0.00%"

# Per-run raw data (rows 5-12) plus the delta formula in column D.
$runData = @(
    @(1, 0.206763,            0.20448),
    @(2, 0.214215999999999,   0.193188),
    @(3, 0.18377,             0.181553999999999),
    @(4, 0.193374,            0.190662999999999),
    @(5, 0.19017,             0.18808),
    @(6, 0.189051,            0.186938999999999),
    @(7, 0.175644,            0.173329999999999),
    @(8, 0.171633,            0.169612)
)

$row = 5
foreach ($r in $runData) {
    $timing.Cells.Item($row, 1).Value = $r[0]
    $timing.Cells.Item($row, 2).Value = $r[1]
    $timing.Cells.Item($row, 3).Value = $r[2]
    # Row 6 is deliberately left without a Delta formula in the source
    # data (matches the original author's incomplete fill-down).
    if ($row -ne 6) {
        $timing.Cells.Item($row, 4).Formula = "=B$row-C$row"
    }
    $row++
}

# Column widths for the two data columns.
$timing.Columns.Item(2).ColumnWidth = 15.2857
$timing.Columns.Item(3).ColumnWidth = 14.2857

# View state: this new sheet is the active / selected one, with B12 as
# the last-selected cell.
$timing.Range("B12").Select()

# ---------------------------------------------------------------------
# 2. Minor leftover view-state tweaks on two pre-existing sheets
#    (selection changes left behind from navigating the workbook).
# ---------------------------------------------------------------------
$orbitalHash = $wb.Worksheets.Item("@orbital_hash")
$orbitalHash.Range("D3:D5").Select()

$ptElementHomeTotal = $wb.Worksheets.Item("custom pt_element_home total")
$ptElementHomeTotal.Range("J3").Select()

# Re-activate the new sheet / cell so it is the one shown & selected
# when the workbook is next opened.
$timing.Activate()
$timing.Range("B12").Select()

$wb.Save()
